# This script applies a weekly data refresh to the "Hortaliza, Femacal de
# La Calera - Sandia" sheet: a brand-new daily record is inserted at the
# top of the existing price-grade block (rows 626-727), which pushes all
# the later grade rows down by one row (row 728 is newly created to hold
# what used to be the last row's data).
#
# Columns D (Fecha) and I..P (Calidad..Precio $/Kg) are the ones that
# belong to each "grade" record and therefore shift down by one row.
# Columns A,B,C,E,F,G,H,Q,R are constant identifying/metadata columns for
# this whole sheet and do not shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 626
$lastOldRow = 727
$lastNewRow = 728

# Column indexes for the columns that shift down by one row.
$colD = 4
$colI = 9
$colJ = 10
$colK = 11
$colL = 12
$colM = 13
$colN = 14
$colO = 15
$colP = 16

# Constant columns (same value on every data row of this sheet).
$colA = 1
$colB = 2
$colC = 3
$colE = 5
$colF = 6
$colG = 7
$colH = 8
$colQ = 17
$colR = 18

# NOTE: this runtime's ".Value" property getter/setter does not behave
# like real Excel COM when chaining a read straight into a write (it
# serializes the property accessor itself instead of the underlying
# value). ".Value2" does not have that problem, so it is used everywhere
# below for both reads and writes.

$constA = $ws.Cells.Item($firstDataRow, $colA).Value2
$constB = $ws.Cells.Item($firstDataRow, $colB).Value2
$constC = $ws.Cells.Item($firstDataRow, $colC).Value2
$constE = $ws.Cells.Item($firstDataRow, $colE).Value2
$constF = $ws.Cells.Item($firstDataRow, $colF).Value2
$constG = $ws.Cells.Item($firstDataRow, $colG).Value2
$constH = $ws.Cells.Item($firstDataRow, $colH).Value2
$constQ = $ws.Cells.Item($firstDataRow, $colQ).Value2
$constR = $ws.Cells.Item($firstDataRow, $colR).Value2

# Make sure the brand new last row (728) has the same date number format
# as the rest of the D column (it doesn't exist yet, so it starts out
# with the generic/default format).
$ws.Cells.Item($lastNewRow, $colD).NumberFormat = $ws.Cells.Item($lastOldRow, $colD).NumberFormat

# Make sure the brand new last row (728) has the constant metadata columns
# populated too, since that row does not exist yet.
$ws.Cells.Item($lastNewRow, $colA).Value2 = $constA
$ws.Cells.Item($lastNewRow, $colB).Value2 = $constB
$ws.Cells.Item($lastNewRow, $colC).Value2 = $constC
$ws.Cells.Item($lastNewRow, $colE).Value2 = $constE
$ws.Cells.Item($lastNewRow, $colF).Value2 = $constF
$ws.Cells.Item($lastNewRow, $colG).Value2 = $constG
$ws.Cells.Item($lastNewRow, $colH).Value2 = $constH
$ws.Cells.Item($lastNewRow, $colQ).Value2 = $constQ
$ws.Cells.Item($lastNewRow, $colR).Value2 = $constR

# Shift the D,I,J,K,L,M,N,O,P values down by one row, working from the
# bottom up so we never overwrite a source row before it has been read.
for ($r = $lastNewRow; $r -gt $firstDataRow; $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, $colD).Value2 = $ws.Cells.Item($src, $colD).Value2
    $ws.Cells.Item($r, $colI).Value2 = $ws.Cells.Item($src, $colI).Value2
    $ws.Cells.Item($r, $colJ).Value2 = $ws.Cells.Item($src, $colJ).Value2
    $ws.Cells.Item($r, $colK).Value2 = $ws.Cells.Item($src, $colK).Value2
    $ws.Cells.Item($r, $colL).Value2 = $ws.Cells.Item($src, $colL).Value2
    $ws.Cells.Item($r, $colM).Value2 = $ws.Cells.Item($src, $colM).Value2
    $ws.Cells.Item($r, $colN).Value2 = $ws.Cells.Item($src, $colN).Value2
    $ws.Cells.Item($r, $colO).Value2 = $ws.Cells.Item($src, $colO).Value2
    $ws.Cells.Item($r, $colP).Value2 = $ws.Cells.Item($src, $colP).Value2
}

# Finally, write the brand new record into row 626.
$ws.Cells.Item($firstDataRow, $colD).Value2 = 45212
$ws.Cells.Item($firstDataRow, $colI).Value2 = "Primera"
$ws.Cells.Item($firstDataRow, $colJ).Value2 = 150
$ws.Cells.Item($firstDataRow, $colK).Value2 = 700
$ws.Cells.Item($firstDataRow, $colL).Value2 = 700
$ws.Cells.Item($firstDataRow, $colM).Value2 = 700
$ws.Cells.Item($firstDataRow, $colN).Value2 = "`$/kilo (volumen en unidades)"
$ws.Cells.Item($firstDataRow, $colO).Value2 = "Perú"
$ws.Cells.Item($firstDataRow, $colP).Value2 = 700
